# Update the benchmark timing results (column B) in the MatrixMultTiming
# sheet with the latest NumPy matrix-multiplication run times.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = [double]"5.483627319335938E-06"
$ws.Range("B3").Value  = [double]"3.862380981445312E-05"
$ws.Range("B4").Value  = [double]"0.0001327991485595703"
$ws.Range("B5").Value  = [double]"0.003339290618896484"
$ws.Range("B6").Value  = [double]"0.01249885559082031"
$ws.Range("B7").Value  = [double]"0.05233240127563477"
$ws.Range("B8").Value  = [double]"0.07263898849487305"
$ws.Range("B9").Value  = [double]"0.3523106575012207"
$ws.Range("B10").Value = [double]"0.9420583248138428"
$ws.Range("B11").Value = [double]"3.11338210105896"
$ws.Range("B12").Value = [double]"7.020748376846313"
